$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking.com crypto price/volume snapshot update (GitHub Actions bot).
# Every D/E cell is stored as text (inline string) in the source sheet, so
# force text format before writing to stop Excel from auto-coercing
# numeric-looking values (e.g. "224.60", "1.00") into Number cells, which
# would silently drop significant trailing zeros.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "34.640.59"
Set-TextValue "D3" "1.790.03"
Set-TextValue "E3" "  +0.50%  "
Set-TextValue "D5" "224.60"
Set-TextValue "E5" "  -0.14%  "
Set-TextValue "D6" "0.560"
Set-TextValue "E6" "  +0.79%  "
Set-TextValue "E7" "  -0.03%  "
Set-TextValue "D8" "32.66"
Set-TextValue "E8" "  +6.66%  "
Set-TextValue "E9" "  +2.48%  "
Set-TextValue "E10" "  +1.37%  "
Set-TextValue "E11" "  +1.42%  "
Set-TextValue "D12" "2.045.52"
Set-TextValue "E12" "  +0.26%  "
Set-TextValue "D13" "11.02"
Set-TextValue "E13" "  +10.71%  "
Set-TextValue "D14" "1.780.84"
Set-TextValue "E14" "  -0.30%  "
Set-TextValue "D15" "0.634"
Set-TextValue "E15" "  +1.14%  "
Set-TextValue "D16" "34.592.36"
Set-TextValue "E16" "  +2.12%  "
Set-TextValue "D17" "4.29"
Set-TextValue "E17" "  +2.65%  "
Set-TextValue "D18" "68.79"
Set-TextValue "E18" "  +0.62%  "
Set-TextValue "D19" "253.97"
Set-TextValue "E19" "  +1.22%  "
Set-TextValue "D20" "0.0₃0767"
Set-TextValue "E20" "  +3.85%  "
Set-TextValue "D21" "1.00"
Set-TextValue "E21" "  -0.03%  "
Set-TextValue "D22" "10.41"
Set-TextValue "E22" "  +1.20%  "
Set-TextValue "D23" "4.24"
Set-TextValue "E23" "  +0.57%  "
Set-TextValue "E24" "  -1.18%  "
Set-TextValue "D25" "159.58"
Set-TextValue "E25" "  +0.27%  "
Set-TextValue "D26" "16.39"
Set-TextValue "E26" "  -0.48%  "
Set-TextValue "D27" "7.09"
Set-TextValue "E27" "  +2.29%  "
Set-TextValue "E28" "  +0.18%  "
Set-TextValue "E29" "  -0.08%  "
Set-TextValue "B30" "Filecoin"
Set-TextValue "C30" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D30" "3.76"
Set-TextValue "E30" "  -0.98%  "
Set-TextValue "B31" "Hedera"
Set-TextValue "C31" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D31" "0.0517"
Set-TextValue "E31" "  +0.66%  "
Set-TextValue "E32" "  +0.28%  "
Set-TextValue "D33" "3.58"
Set-TextValue "E33" "  +0.89%  "
Set-TextValue "E34" "  +3.23%  "
Set-TextValue "D35" "1.444.56"
Set-TextValue "E35" "  -2.59%  "
Set-TextValue "E36" "  -0.17%  "
Set-TextValue "E37" "  +2.61%  "
Set-TextValue "E38" "  -0.57%  "
Set-TextValue "D39" "83.08"
Set-TextValue "E39" "  -0.31%  "
Set-TextValue "E40" "  +3.97%  "
Set-TextValue "E41" "  -0.19%  "
Set-TextValue "D42" "0.900"
Set-TextValue "E42" "  +1.52%  "
Set-TextValue "E43" "  -0.51%  "
Set-TextValue "E44" "  -0.74%  "
Set-TextValue "B45" "FraxShare"
Set-TextValue "C45" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D45" "5.90"
Set-TextValue "E45" "  +2.28%  "
Set-TextValue "B46" "WEMIXToken"
Set-TextValue "C46" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D46" "1.05"
Set-TextValue "E46" "  -1.60%  "
Set-TextValue "D47" "1.941.05"
Set-TextValue "E47" "  +0.08%  "
Set-TextValue "B48" "InjectiveProtocol"
Set-TextValue "C48" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D48" "12.02"
Set-TextValue "E48" "  +1.00%  "
Set-TextValue "B49" "PaxDollar"
Set-TextValue "C49" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D49" "1.00"
Set-TextValue "E49" "  -0.04%  "
Set-TextValue "D50" "103.25"
Set-TextValue "E50" "  +5.81%  "
Set-TextValue "E51" "  +4.94%  "
